# Add a new Job Posting row (Job_Id = 16) to the LinkedIn job postings sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the current data block.
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = 16
$ws.Cells.Item($newRow, 2).Value = "stack Developer"
$ws.Cells.Item($newRow, 3).Value = "gggggghfdgf"
$ws.Cells.Item($newRow, 4).Value = 1
$ws.Cells.Item($newRow, 5).Value = 2
$ws.Cells.Item($newRow, 6).Value = 0
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0
